# Sync automatico del tracker (cada 3h)
# Appends the latest batch of tracked match rows to the bottom of the
# results tracker sheet, following the same column layout as the
# existing data: event_id, fecha, jugador_A, jugador_B, pronostico,
# cuota, resultado, profit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,@(361, 14428732, '2025-08-21', 'Sebastian Korda', 'Miomir Kecmanovic', 'Gana Miomir Kecmanovic', 2.1)
    ,@(362, 14428725, '2025-08-21', 'Yunchaokete Bu', 'Botic Van de Zandschulp', 'Gana Botic Van de Zandschulp', 2)
    ,@(363, 14427813, '2025-08-21', 'Diana Shnaider', 'Elise Mertens', 'Gana Diana Shnaider', 2.1)
    ,@(364, 14427818, '2025-08-21', 'Alycia Parks', 'Rebecca Sramkova', 'Gana Alycia Parks', 2.1)
    ,@(365, 14427815, '2025-08-22', 'Linda Noskova', 'Ekaterina Alexandrova', 'Gana Linda Noskova', 2.63)
    ,@(366, 14486673, '2025-08-21', 'Martin Damm Jr', 'Benjamin Hassan', 'Gana Benjamin Hassan', 1.4)
    ,@(367, 14485932, '2025-08-21', 'Beibit Zhukayev', 'Zachary Svajda', 'Gana Beibit Zhukayev', 3.5)
    ,@(368, 14486400, '2025-08-21', 'Francesca Jones', 'Ekaterine Gorgodze', 'Gana Ekaterine Gorgodze', 4.5)
    ,@(369, 14485925, '2025-08-21', 'Cadence Brace', 'Veronika Erjavec', 'Gana Cadence Brace', 2.62)
    ,@(370, 14486278, '2025-08-21', 'Carol Zhao', 'Lucrezia Stefanini', 'Gana Carol Zhao', 2.75)
    ,@(371, 14477278, '2025-08-21', 'Claire Liu', 'Maddison Inglis', 'Gana Maddison Inglis', 1.62)
    ,@(372, 14418940, '2025-08-21', 'Jakub Paul', 'Robin Bertrand', 'Gana Robin Bertrand', 2.63)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]            # event_id
    $ws.Cells.Item($r, 2).Value = "'" + $row[2]       # fecha (kept as text, like existing rows)
    $ws.Cells.Item($r, 3).Value = $row[3]             # jugador_A
    $ws.Cells.Item($r, 4).Value = $row[4]             # jugador_B
    $ws.Cells.Item($r, 5).Value = $row[5]             # pronostico
    $ws.Cells.Item($r, 6).Value = $row[6]             # cuota
    $ws.Cells.Item($r, 7).Value = ""                  # resultado (pending)
    $ws.Cells.Item($r, 8).Value = ""                  # profit (pending)
}
